$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Delete the "Meta description: Discover Beetle Jewels..." paragraph that
#    currently sits right after the H1 title at the top of the document.
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
if ($metaPara.Range.Text -like "Meta description*") {
    $metaPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2. Insert a new bold paragraph "Play Beetle Jewels Free: Review of RTP,
#    Bonus Features & Graphics" right before the final paragraph (the one
#    that currently holds the "Prompt for Feature Image..." text).
#
#    We use Range.InsertXML so we get exact control over the run layout
#    (a leading empty run followed by a single bold run), matching the
#    structure used elsewhere in this document.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$insertPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$xmlFragment = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Play Beetle Jewels Free: Review of RTP, Bonus Features &amp; Graphics</w:t></w:r></w:p><w:p/><w:sectPr><w:pgSz w:w="12240" w:h="15840"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$insertPoint.InsertXML($xmlFragment)

# InsertXML leaves behind one extra empty paragraph mark after our new
# paragraph (an artifact of inserting a full <w:p> in the middle of text) -
# remove it so paragraph numbering/structure is back to normal. (Note: an
# "empty" paragraph's Range.Text is the single paragraph-mark character, not
# a zero-length string, so check Length <= 1 rather than equality with "".)
$strayPara = $d.Paragraphs.Item($count + 1)
if ($strayPara.Range.Text.Length -le 1) {
    $strayPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 3. Replace the text of the final paragraph (formerly "Prompt for Feature
#    Image: ...") with the meta-description copy, keeping its italic run
#    formatting untouched.
# ---------------------------------------------------------------------------
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalPara.Range.Find.Execute(
    'Prompt for Feature Image: Create a fun and colorful cartoon style image that features a happy Maya warrior wearing glasses and holding a precious beetle jewel. The warrior should be surrounded by a lush forest filled with colorful beetles. The image should reflect the enchanted and mystical theme of the game, while also highlighting the beetle symbols and the importance of the precious stones they carry. Use vibrant colors and playful details to attract players and make the image stand out. Don''t forget to include the game title, "Beetle Jewels," in a bold and eye-catching font that complements the overall design.',
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discover Beetle Jewels, a slot machine game with exciting Wild and Scatter symbols, high RTP, and striking graphics. Play now for free!",
    2)
